$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data rows 2-4 (participants 1-3) below the header row 1.
# Extend the participant list with 6 more blank rows (rows 5-10), each with a
# sequential id in column A (4..9) and empty (but present) placeholder cells in
# columns B..P, mirroring the existing blank participant rows 3-4.
$lastCol = 16  # column P

for ($rowIndex = 5; $rowIndex -le 10; $rowIndex++) {
    $participantId = $rowIndex - 1
    $ws.Cells.Item($rowIndex, 1).Value = $participantId

    for ($col = 2; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        # A bare empty-string assignment is treated as "no value" and the
        # cell would not be written out at all. Seeding it with a leading
        # apostrophe forces the engine to materialise a real (blank) text
        # cell - same as the pre-existing placeholder cells - then we reset
        # the style back to Normal so no stray cell formatting is left
        # behind.
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
